# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# Only the Price (D) and Volume(1h) (E) columns change; Coin/Link/rank stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.540.66"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "2.469.95"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'321.71"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").Value = "'105.18"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("D7").Value = "'0.522"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").Value = "'36.22"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'18.30"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "2.862.94"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "2.476.02"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "46.446.53"
$ws.Range("E18").Value = "  +4.23%  "
$ws.Range("D19").Value = "'12.68"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("D22").Value = "'70.64"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").Value = "'248.94"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").Value = "'26.15"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").Value = "'34.74"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").Value = "'49.69"
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "'19.68"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'5.33"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "'4.63"
$ws.Range("E37").Value = "  +3.31%  "
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").Value = "'123.64"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "'20.77"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "1.985.50"
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "'2.09"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E48").Value = "  +7.24%  "
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").Value = "'5.26"
$ws.Range("E50").Value = "  +13.31%  "
$ws.Range("D51").Value = "'79.05"
$ws.Range("E51").Value = "  +5.18%  "
